$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
try {
    $newShp = $ftr.Range.InlineShapes.AddPicture("C:\fake\image2.png")
    Write-Output "Added: $($newShp.AlternativeText)"
} catch {
    Write-Output "ERROR: $_"
}
